$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number must be forced to Text
# so Excel does not silently convert them to a numeric type (matches the
# original inlineStr/text storage in the workbook).
function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '65.403.52'
$ws.Range('E2').Value = '  -3.60%  '
$ws.Range('D3').Value = '3.487.02'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('E4').Value = '  +0.06%  '
Set-TextCell 'D5' '552.02'
$ws.Range('E5').Value = '  -0.92%  '
Set-TextCell 'D6' '178.98'
$ws.Range('E6').Value = '  -6.77%  '
Set-TextCell 'D7' '0.640'
$ws.Range('E7').Value = '  +4.47%  '
$ws.Range('E8').Value = '  +0.04%  '
Set-TextCell 'D9' '0.631'
$ws.Range('E9').Value = '  -1.37%  '
$ws.Range('E10').Value = '  +1.86%  '
Set-TextCell 'D11' '53.68'
$ws.Range('E11').Value = '  -6.14%  '
$ws.Range('E12').Value = '  -2.42%  '
Set-TextCell 'D13' '9.17'
$ws.Range('E13').Value = '  -3.50%  '
$ws.Range('D14').Value = '4.043.57'
$ws.Range('E14').Value = '  -0.79%  '
$ws.Range('D15').Value = '3.488.79'
$ws.Range('E15').Value = '  -0.98%  '
$ws.Range('E16').Value = '  +0.04%  '
Set-TextCell 'D17' '18.36'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('E18').Value = '  +1.78%  '
$ws.Range('D19').Value = '65.450.88'
$ws.Range('E19').Value = '  -4.60%  '
Set-TextCell 'D20' '0.992'
$ws.Range('E20').Value = '  -1.81%  '
Set-TextCell 'D21' '412.58'
$ws.Range('E21').Value = '  +0.89%  '
Set-TextCell 'D22' '4.03'
$ws.Range('E22').Value = '  +1.73%  '
Set-TextCell 'D23' '85.65'
$ws.Range('E23').Value = '  +0.85%  '
Set-TextCell 'D24' '4.09'
$ws.Range('E24').Value = '  -3.31%  '
$ws.Range('E25').Value = '  +6.03%  '
Set-TextCell 'D26' '10.76'
$ws.Range('E26').Value = '  -7.93%  '
Set-TextCell 'D27' '2.84'
$ws.Range('E27').Value = '  -2.27%  '
$ws.Range('E28').Value = '  -1.88%  '
Set-TextCell 'D29' '9.01'
$ws.Range('E29').Value = '  +4.11%  '
Set-TextCell 'D30' '30.20'
$ws.Range('E30').Value = '  -1.31%  '
Set-TextCell 'D31' '612.98'
$ws.Range('E31').Value = '  -10.47%  '
Set-TextCell 'D32' '6.44'
$ws.Range('E32').Value = '  -6.50%  '
Set-TextCell 'D33' '11.63'
$ws.Range('E33').Value = '  -0.92%  '
$ws.Range('E34').Value = '  -1.58%  '
Set-TextCell 'D35' '59.37'
$ws.Range('E35').Value = '  -1.89%  '
$ws.Range('E36').Value = '  +11.19%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').Value = '  -5.53%  '
$ws.Range('D39').Value = '0.0₃0786'
$ws.Range('E39').Value = '  -7.10%  '
$ws.Range('D40').Value = '3.361.11'
$ws.Range('E40').Value = '  +10.18%  '
$ws.Range('E41').Value = '  -6.33%  '
Set-TextCell 'D42' '1.00'
$ws.Range('E42').Value = '  +0.14%  '
Set-TextCell 'D43' '3.26'
$ws.Range('E43').Value = '  -4.95%  '
$ws.Range('E44').Value = '  -6.54%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell 'D45' '2.51'
$ws.Range('E45').Value = '  -9.10%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell 'D46' '3.25'
$ws.Range('E46').Value = '  +2.27%  '
$ws.Range('E47').Value = '  -2.12%  '
Set-TextCell 'D48' '2.75'
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('E49').Value = '  +1.59%  '
Set-TextCell 'D50' '137.52'
$ws.Range('E50').Value = '  -1.03%  '
$ws.Range('E51').Value = '  -9.97%  '
